$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Exam Duration" header to "Exam Duration (Hours)"
$ws.Range("I1").Value = "Exam Duration (Hours)"

# Widen column I to fit the new, longer header text
$ws.Columns("I").ColumnWidth = 22.17

# Leave the selection where the author last left it
[void]$ws.Range("N8").Select()
